# Auto-generated edit script: update FFXIV leve market-price figures (columns H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets, per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Cells.Item(17, 8).Value = 1564547.6
$ws.Cells.Item(17, 10).Value = 1589371.2
$ws.Cells.Item(17, 12).Value = 4768113.6
$ws.Cells.Item(17, 14).Value = -4768449.6

# Row 58
$ws.Cells.Item(58, 8).Value = 3236.8462
$ws.Cells.Item(58, 9).Value = 309.0
$ws.Cells.Item(58, 10).Value = 5066.75
$ws.Cells.Item(58, 11).Value = 927.0
$ws.Cells.Item(58, 12).Value = 15200.25
$ws.Cells.Item(58, 13).Value = -777.0
$ws.Cells.Item(58, 14).Value = -15500.25

# Row 88
$ws.Cells.Item(88, 8).Value = 1314.0526
$ws.Cells.Item(88, 9).Value = 1034.5714
$ws.Cells.Item(88, 10).Value = 1477.0834
$ws.Cells.Item(88, 11).Value = 1034.5714
$ws.Cells.Item(88, 12).Value = 1477.0834
$ws.Cells.Item(88, 13).Value = -628.5714
$ws.Cells.Item(88, 14).Value = -2289.0834

# Row 91
$ws.Cells.Item(91, 8).Value = 1314.0526
$ws.Cells.Item(91, 9).Value = 1034.5714
$ws.Cells.Item(91, 10).Value = 1477.0834
$ws.Cells.Item(91, 11).Value = 1034.5714
$ws.Cells.Item(91, 12).Value = 1477.0834
$ws.Cells.Item(91, 13).Value = 369.4286
$ws.Cells.Item(91, 14).Value = -4285.0834

# Row 116
$ws.Cells.Item(116, 8).Value = 3319.125
$ws.Cells.Item(116, 9).Value = 1432.5
$ws.Cells.Item(116, 10).Value = 4451.1
$ws.Cells.Item(116, 11).Value = 1432.5
$ws.Cells.Item(116, 12).Value = 4451.1
$ws.Cells.Item(116, 13).Value = 2009.5
$ws.Cells.Item(116, 14).Value = -11335.1

# Row 129
$ws.Cells.Item(129, 8).Value = 263970.16
$ws.Cells.Item(129, 9).Value = 400.0
$ws.Cells.Item(129, 10).Value = 278612.94
$ws.Cells.Item(129, 11).Value = 1200.0
$ws.Cells.Item(129, 12).Value = 835838.8200000001
$ws.Cells.Item(129, 13).Value = 3800.0
$ws.Cells.Item(129, 14).Value = -845838.8200000001

# Row 135
$ws.Cells.Item(135, 8).Value = 20841830.0
$ws.Cells.Item(135, 9).Value = 1056.1875
$ws.Cells.Item(135, 11).Value = 9505.6875
$ws.Cells.Item(135, 13).Value = -6970.6875

# Row 138
$ws.Cells.Item(138, 8).Value = 1430.8889
$ws.Cells.Item(138, 9).Value = 510.54544
$ws.Cells.Item(138, 10).Value = 2877.1428
$ws.Cells.Item(138, 11).Value = 1531.63632
$ws.Cells.Item(138, 12).Value = 8631.4284
$ws.Cells.Item(138, 13).Value = 3608.36368
$ws.Cells.Item(138, 14).Value = -18911.4284

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Cells.Item(45, 8).Value = 2375.1428
$ws.Cells.Item(45, 9).Value = 2435.6924
$ws.Cells.Item(45, 11).Value = 2435.6924
$ws.Cells.Item(45, 13).Value = -2058.6924

# Row 132
$ws.Cells.Item(132, 8).Value = 12589.17
$ws.Cells.Item(132, 9).Value = 1689.1212
$ws.Cells.Item(132, 10).Value = 38282.145
$ws.Cells.Item(132, 11).Value = 5067.363600000001
$ws.Cells.Item(132, 12).Value = 114846.435
$ws.Cells.Item(132, 13).Value = -2537.363600000001
$ws.Cells.Item(132, 14).Value = -119906.435

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Cells.Item(99, 8).Value = 1582.7142
$ws.Cells.Item(99, 9).Value = 1318.625
$ws.Cells.Item(99, 11).Value = 1318.625
$ws.Cells.Item(99, 13).Value = 179.375

# Row 132
$ws.Cells.Item(132, 8).Value = 0.0
$ws.Cells.Item(132, 10).Value = 0.0
$ws.Cells.Item(132, 14).ClearContents()

# Row 134
$ws.Cells.Item(134, 8).Value = 32281.572
$ws.Cells.Item(134, 9).Value = 39840.285
$ws.Cells.Item(134, 10).Value = 2046.7142
$ws.Cells.Item(134, 11).Value = 119520.855
$ws.Cells.Item(134, 12).Value = 6140.142599999999
$ws.Cells.Item(134, 13).Value = -116985.855
$ws.Cells.Item(134, 14).Value = -11210.1426

$ws = $wb.Worksheets.Item("CRP")
# Row 28
$ws.Cells.Item(28, 8).Value = 0.0
$ws.Cells.Item(28, 10).Value = 0.0
$ws.Cells.Item(28, 14).ClearContents()

# Row 94
$ws.Cells.Item(94, 8).Value = 2468.9524
$ws.Cells.Item(94, 9).Value = 1469.1666
$ws.Cells.Item(94, 11).Value = 1469.1666
$ws.Cells.Item(94, 13).Value = -1018.1666

# Row 99
$ws.Cells.Item(99, 8).Value = 33338914.0
$ws.Cells.Item(99, 9).Value = 4462.5
$ws.Cells.Item(99, 10).Value = 71435430.0
$ws.Cells.Item(99, 11).Value = 4462.5
$ws.Cells.Item(99, 12).Value = 71435430.0
$ws.Cells.Item(99, 13).Value = -2964.5
$ws.Cells.Item(99, 14).Value = -71438426.0

# Row 126
$ws.Cells.Item(126, 8).Value = 33338914.0
$ws.Cells.Item(126, 9).Value = 4462.5
$ws.Cells.Item(126, 10).Value = 71435430.0
$ws.Cells.Item(126, 11).Value = 13387.5
$ws.Cells.Item(126, 12).Value = 214306290.0
$ws.Cells.Item(126, 13).Value = -10917.5
$ws.Cells.Item(126, 14).Value = -214311230.0

# Row 134
$ws.Cells.Item(134, 8).Value = 1190.6938
$ws.Cells.Item(134, 9).Value = 900.68
$ws.Cells.Item(134, 10).Value = 1492.7916
$ws.Cells.Item(134, 11).Value = 2702.04
$ws.Cells.Item(134, 12).Value = 4478.3748
$ws.Cells.Item(134, 13).Value = -167.04
$ws.Cells.Item(134, 14).Value = -9548.3748

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 1552.5
$ws.Cells.Item(5, 9).Value = 1077.5
$ws.Cells.Item(5, 10).Value = 2502.5
$ws.Cells.Item(5, 11).Value = 3232.5
$ws.Cells.Item(5, 12).Value = 7507.5
$ws.Cells.Item(5, 13).Value = -3120.5
$ws.Cells.Item(5, 14).Value = -7731.5

# Row 70
$ws.Cells.Item(70, 8).Value = 4288.706
$ws.Cells.Item(70, 9).Value = 3424.889
$ws.Cells.Item(70, 10).Value = 5260.5
$ws.Cells.Item(70, 11).Value = 10274.667
$ws.Cells.Item(70, 12).Value = 15781.5
$ws.Cells.Item(70, 13).Value = -9959.667000000001
$ws.Cells.Item(70, 14).Value = -16411.5

# Row 73
$ws.Cells.Item(73, 8).Value = 4288.706
$ws.Cells.Item(73, 9).Value = 3424.889
$ws.Cells.Item(73, 10).Value = 5260.5
$ws.Cells.Item(73, 11).Value = 10274.667
$ws.Cells.Item(73, 12).Value = 15781.5
$ws.Cells.Item(73, 13).Value = -9182.667000000001
$ws.Cells.Item(73, 14).Value = -17965.5

# Row 112
$ws.Cells.Item(112, 8).Value = 1545.4
$ws.Cells.Item(112, 9).Value = 931.75
$ws.Cells.Item(112, 11).Value = 2795.25
$ws.Cells.Item(112, 13).Value = -1687.25

# Row 122
$ws.Cells.Item(122, 8).Value = 810.5833
$ws.Cells.Item(122, 10).Value = 1009.625
$ws.Cells.Item(122, 12).Value = 9086.625
$ws.Cells.Item(122, 14).Value = -13986.625

# Row 129
$ws.Cells.Item(129, 8).Value = 294721.53
$ws.Cells.Item(129, 9).Value = 497.14285
$ws.Cells.Item(129, 10).Value = 500678.6
$ws.Cells.Item(129, 11).Value = 1491.42855
$ws.Cells.Item(129, 12).Value = 1502035.8
$ws.Cells.Item(129, 13).Value = 3508.57145
$ws.Cells.Item(129, 14).Value = -1512035.8

# Row 131
$ws.Cells.Item(131, 8).Value = 747.78
$ws.Cells.Item(131, 10).Value = 758.7263
$ws.Cells.Item(131, 12).Value = 2276.1789
$ws.Cells.Item(131, 14).Value = -12356.1789

# Row 135
$ws.Cells.Item(135, 8).Value = 1552.5
$ws.Cells.Item(135, 9).Value = 1077.5
$ws.Cells.Item(135, 10).Value = 2502.5
$ws.Cells.Item(135, 11).Value = 9697.5
$ws.Cells.Item(135, 12).Value = 22522.5
$ws.Cells.Item(135, 13).Value = -7162.5
$ws.Cells.Item(135, 14).Value = -27592.5

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Cells.Item(80, 8).Value = 3355.7273
$ws.Cells.Item(80, 9).Value = 2934.6667
$ws.Cells.Item(80, 10).Value = 3647.2307
$ws.Cells.Item(80, 11).Value = 2934.6667
$ws.Cells.Item(80, 12).Value = 3647.2307
$ws.Cells.Item(80, 13).Value = -1936.6667
$ws.Cells.Item(80, 14).Value = -5643.2307

# Row 83
$ws.Cells.Item(83, 8).Value = 3355.7273
$ws.Cells.Item(83, 9).Value = 2934.6667
$ws.Cells.Item(83, 10).Value = 3647.2307
$ws.Cells.Item(83, 11).Value = 14673.3335
$ws.Cells.Item(83, 12).Value = 18236.1535
$ws.Cells.Item(83, 13).Value = -9681.3335
$ws.Cells.Item(83, 14).Value = -28220.1535

# Row 132
$ws.Cells.Item(132, 8).Value = 38779.285
$ws.Cells.Item(132, 9).Value = 39084.43
$ws.Cells.Item(132, 10).Value = 38169.0
$ws.Cells.Item(132, 11).Value = 117253.29
$ws.Cells.Item(132, 12).Value = 114507.0
$ws.Cells.Item(132, 13).Value = -114723.29
$ws.Cells.Item(132, 14).Value = -119567.0

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Cells.Item(61, 9).Value = 3512.5
$ws.Cells.Item(61, 10).Value = 8400.6
$ws.Cells.Item(61, 11).Value = 3512.5
$ws.Cells.Item(61, 12).Value = 8400.6
$ws.Cells.Item(61, 13).Value = -3310.5
$ws.Cells.Item(61, 14).Value = -8804.6

# Row 113
$ws.Cells.Item(113, 9).Value = 3512.5
$ws.Cells.Item(113, 10).Value = 8400.6
$ws.Cells.Item(113, 11).Value = 3512.5
$ws.Cells.Item(113, 12).Value = 8400.6
$ws.Cells.Item(113, 13).Value = -1342.5
$ws.Cells.Item(113, 14).Value = -12740.6
